$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2026-01-31 05:10"
$ws.Range("B11").Value = 23
$ws.Range("C11").Value = 5
